$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header summary figures (Valor Mora total / Cant. Trabajadores / Cant.
#    Periodos) bumped to reflect the new worker being added to the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 306653
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 6

# ---------------------------------------------------------------------------
# 2) Make room for a new data row. The table currently ends at row 21 (the
#    last row carries a distinct "closing" border style). Insert a fresh row
#    at 22 so everything below (the signature block) shifts down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).Insert()

# Row 22 (the new last row) should get the special "closing" border style
# that row 21 currently has - copy it (values + formatting) down first.
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))

# Row 21 is no longer the last row of the table, so it should switch back to
# the regular interior row style - copy that down from row 20.
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# ---------------------------------------------------------------------------
# 3) Rewrite the data rows (16-22) with the reordered / updated records.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "33104272"
$ws.Range("D16").Value = "NISIDA ALCAZAR JIMENEZ"
$ws.Range("E16").Value = "2309"
$ws.Range("F16").Value = 3627
$ws.Range("G16").Value = 1360000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1043962314"
$ws.Range("D17").Value = "BANIC HAISAWA OTERO LABRADOR"
$ws.Range("E17").Value = "2309"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1160000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73581883"
$ws.Range("D18").Value = "JORGE LUIS MARTINEZ ORTIZ"
$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 38666
$ws.Range("G18").Value = 1000000

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "4865522"
$ws.Range("D19").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E19").Value = "2303"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1680000

$ws.Range("B20").Value = "PPT"
$ws.Range("C20").Value = "4865522"
$ws.Range("D20").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 60000
$ws.Range("G20").Value = 1680000

$ws.Range("B21").Value = "PPT"
$ws.Range("C21").Value = "4865522"
$ws.Range("D21").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E21").Value = "2301"
$ws.Range("F21").Value = 60000
$ws.Range("G21").Value = 1680000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1085045171"
$ws.Range("D22").Value = "LUIS ROBERTO RAMIREZ NAVARRO"
$ws.Range("E22").Value = "2502"
$ws.Range("F22").Value = 37960
$ws.Range("G22").Value = 1423500

# Clear any stray H/I/J leftovers on the data rows (Observaciones / Novedad
# columns are left blank for every worker row in this report).
$ws.Range("H16:J22").ClearContents()
